$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:29 down to 14:30.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new weekly record.
$ws.Range("A13").Value = 3
$ws.Range("B13").Value = "Femacal de La Calera"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44483
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 100112022
$ws.Range("G13").Value = "Arveja Verde"
$ws.Range("H13").Value = "Perfection"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 55
$ws.Range("K13").Value = 29000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = 29455
$ws.Range("N13").Value = "`$/malla 25 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 1178
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
